# "Fix addFonted and createNewSheet"
# Adds a second worksheet ("File too") that summarises only the TOP-level
# tasks from "Tasks" (rows where column B == "TOP"), dropping the Level
# column, and re-drawing the "x" markers as blank/space cells styled with
# a dedicated font + white fill + left alignment ("addFonted" cell style).

$wb = $excel.ActiveWorkbook
$tasks = $wb.Worksheets.Item("Tasks")

# ---------------------------------------------------------------------
# 1. Create the new sheet right after "Tasks" and name it.
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $tasks)
$newSheet.Name = "File too"

# ---------------------------------------------------------------------
# 2. Header row (row 1): same 6 text headers as "Tasks" plus day numbers
#    1..20 across columns G:Z.
# ---------------------------------------------------------------------
$headers = @("      ", "Level", "Id", "Description", "Cost", "Effort")
for ($c = 1; $c -le 6; $c++) {
    $newSheet.Cells.Item(1, $c).Value = $headers[$c - 1]
}
for ($d = 1; $d -le 20; $d++) {
    $newSheet.Cells.Item(1, 6 + $d).Value = $d
}

# ---------------------------------------------------------------------
# 3. Copy over only the TOP level rows from "Tasks" (column B = "TOP"),
#    dropping column B (Level) -- C:F shift into C:F, and re-create the
#    "x" markers (now a single space) in the same G:Z columns they had.
#    Styling (font/fill/alignment) is applied per marker cell so empty
#    cells in between stay untouched/unwritten.
# ---------------------------------------------------------------------
$srcRows = 2, 9, 12
$dstRow = 2
foreach ($srcRow in $srcRows) {
    for ($c = 3; $c -le 6; $c++) {
        $newSheet.Cells.Item($dstRow, $c).Value = $tasks.Cells.Item($srcRow, $c).Value()
    }
    for ($c = 7; $c -le 26; $c++) {
        $v = $tasks.Cells.Item($srcRow, $c).Value()
        if ($v -ne $null) {
            $mark = $newSheet.Cells.Item($dstRow, $c)
            $mark.Value = " "
            $mark.Interior.Pattern = 1
            $mark.Interior.ColorIndex = 2
            $mark.Font.Name = "Aptos Narrow"
            $mark.HorizontalAlignment = -4131
            $mark.WrapText = $false
        }
    }
    $dstRow = $dstRow + 1
}

# ---------------------------------------------------------------------
# 4. Style the header row: new font "Aptos Narrow" + solid white fill +
#    left align + no wrap.
# ---------------------------------------------------------------------
$styleRange = $newSheet.Range($newSheet.Cells.Item(1, 1), $newSheet.Cells.Item(1, 26))
$styleRange.Interior.Pattern = 1
$styleRange.Interior.ColorIndex = 2
$styleRange.Font.Name = "Aptos Narrow"
$styleRange.HorizontalAlignment = -4131
$styleRange.WrapText = $false

# ---------------------------------------------------------------------
# 5. Re-fit columns on both sheets and restore "Tasks" as the active tab.
# ---------------------------------------------------------------------
$newSheet.Columns.AutoFit()
$tasks.Columns.AutoFit()
$tasks.Activate()
